$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.540.22"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.443.36"
$ws.Range("E3").Value = "  -2.78%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.26"
$ws.Range("E5").Value = "  -2.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.07"
$ws.Range("E6").Value = "  -3.62%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.441.22"
$ws.Range("E9").Value = "  -2.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -2.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.96"
$ws.Range("E11").Value = "  -2.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -4.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.044.63"
$ws.Range("E13").Value = "  -2.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.84"
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.526.44"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000175"
$ws.Range("E17").Value = "  -4.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.444.85"
$ws.Range("E18").Value = "  -2.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("E19").Value = "  -4.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.91"
$ws.Range("E20").Value = "  -7.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.32"
$ws.Range("E21").Value = "  -3.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.80"
$ws.Range("E22").Value = "  -4.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +1.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.25"
$ws.Range("E25").Value = "  -3.60%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.530"
$ws.Range("E26").Value = "  -3.41%  "

$ws.Range("E27").Value = "  -5.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  -5.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.174"
$ws.Range("E29").Value = "  -3.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.03"
$ws.Range("E31").Value = "  -4.76%  "

$ws.Range("E32").Value = "  -2.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  -6.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.31"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.20"
$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("E37").Value = "  -8.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.22"
$ws.Range("E38").Value = "  -2.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.880"
$ws.Range("E39").Value = "  -0.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  -5.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").Value = "  -4.39%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.60"
$ws.Range("E42").Value = "  -8.21%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.51"
$ws.Range("E43").Value = "  -5.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.77"
$ws.Range("E44").Value = "  -6.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0710"
$ws.Range("E45").Value = "  -4.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.87"
$ws.Range("E46").Value = "  -7.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.692.22"
$ws.Range("E47").Value = "  -7.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.14"
$ws.Range("E48").Value = "  -3.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0295"
$ws.Range("E49").Value = "  -4.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "322.35"
$ws.Range("E50").Value = "  -9.06%  "

$ws.Range("E51").Value = "  -5.94%  "
